# Insert a new weekly price record for "Poroto granado" (Macroferia Regional
# de Talca) right before the current row 109. Excel shifts the existing
# rows 109-132 down to 110-133 automatically, preserving their contents and
# formatting, and the sheet's used-range grows to A1:R133.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(109).Insert()

# Populate the newly inserted (blank) row 109 with the new record.
$ws.Range("A109").Value = 5
$ws.Range("B109").Value = "Macroferia Regional de Talca"
$ws.Range("C109").Value = "Maule"
$ws.Range("D109").Value = 44637
$ws.Range("E109").Value = 7
$ws.Range("F109").Value = 100112030
$ws.Range("G109").Value = "Poroto granado"
$ws.Range("H109").Value = "Sin especificar"
$ws.Range("I109").Value = "Primera"
$ws.Range("J109").Value = 300
$ws.Range("K109").Value = 18000
$ws.Range("L109").Value = 18000
$ws.Range("M109").Value = 18000
$ws.Range("N109").Value = "`$/saco 25 kilos"
$ws.Range("O109").Value = "Región del Maule"
$ws.Range("P109").Value = 720
$ws.Range("Q109").Value = 25
$ws.Range("R109").Value = "Hortaliza"
